$wb = $excel.ActiveWorkbook

# --- Sheets "Item"/"Activity": rename attribute & reward header columns
#     (finish conversion of list). The write order below matches the
#     order new labels were appended to the shared-string table in the
#     authored workbook.
$itemWs = $wb.Worksheets.Item("Item")
$activityWs = $wb.Worksheets.Item("Activity")

$itemWs.Range("E1").Value = "属性1ID"
$itemWs.Range("H1").Value = "属性2Value"
$itemWs.Range("G1").Value = "属性2ID"
$activityWs.Range("E1").Value = "奖励1ID"
$activityWs.Range("H1").Value = "奖励2NUM"
$activityWs.Range("G1").Value = "奖励2ID"
$activityWs.Range("F1").Value = "奖励1NUM"
$itemWs.Range("F1").Value = "属性1Value"

# --- View/selection bookkeeping to mirror the saved workbook state ---
$activityWs.Activate()
$activityWs.Range("G12").Select()

$itemWs.Activate()
$itemWs.Range("J13").Select()
